$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1050
$ws.Range("I18").Value = 1050
$ws.Range("K18").Value = 1050
$ws.Range("M18").Value = -766

$ws.Range("H74").Value = 10499.75
$ws.Range("I74").Value = 7999
$ws.Range("J74").Value = 11333.333
$ws.Range("K74").Value = 7999
$ws.Range("L74").Value = 11333.333
$ws.Range("M74").Value = -7063
$ws.Range("N74").Value = -13205.333

$ws.Range("H77").Value = 10499.75
$ws.Range("I77").Value = 7999
$ws.Range("J77").Value = 11333.333
$ws.Range("K77").Value = 39995
$ws.Range("L77").Value = 56666.665
$ws.Range("M77").Value = -35315
$ws.Range("N77").Value = -66026.66500000001

$ws.Range("H116").Value = 3857.6667
$ws.Range("I116").Value = 3828.5715
$ws.Range("J116").Value = 3898.4
$ws.Range("K116").Value = 3828.5715
$ws.Range("L116").Value = 3898.4
$ws.Range("M116").Value = -386.5715
$ws.Range("N116").Value = -10782.4

$ws.Range("H132").Value = 36532.668
$ws.Range("I132").Value = 36532.668
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 109598.004
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -107068.004
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 99581.664
$ws.Range("J133").Value = 99581.664
$ws.Range("L133").Value = 99581.664
$ws.Range("N133").Value = -109701.664

$ws.Range("H137").Value = 1042.75
$ws.Range("I137").Value = 835.58826
$ws.Range("K137").Value = 2506.76478
$ws.Range("M137").Value = 43.23522000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3965.8462
$ws.Range("I2").Value = 999.3333
$ws.Range("J2").Value = 6508.5713
$ws.Range("K2").Value = 999.3333
$ws.Range("L2").Value = 6508.5713
$ws.Range("M2").Value = -886.3333
$ws.Range("N2").Value = -6734.5713

$ws.Range("H32").Value = 3666.6287
$ws.Range("I32").Value = 3576.7424
$ws.Range("K32").Value = 3576.7424
$ws.Range("M32").Value = -3289.7424

$ws.Range("H45").Value = 6050.5386
$ws.Range("I45").Value = 3644.6667
$ws.Range("K45").Value = 3644.6667
$ws.Range("M45").Value = -3267.6667

$ws.Range("H46").Value = 7139.5454
$ws.Range("J46").Value = 7403.8887
$ws.Range("L46").Value = 7403.8887
$ws.Range("N46").Value = -8041.8887

$ws.Range("H61").Value = 13164247
$ws.Range("I61").Value = 15156668
$ws.Range("K61").Value = 15156668
$ws.Range("M61").Value = -15156456

$ws.Range("H116").Value = 3965.8462
$ws.Range("I116").Value = 999.3333
$ws.Range("J116").Value = 6508.5713
$ws.Range("K116").Value = 999.3333
$ws.Range("L116").Value = 6508.5713
$ws.Range("M116").Value = 1294.6667
$ws.Range("N116").Value = -11096.5713

$ws.Range("H130").Value = 85830
$ws.Range("J130").Value = 85830
$ws.Range("L130").Value = 85830
$ws.Range("N130").Value = -95870

$ws.Range("H132").Value = 3820.721
$ws.Range("I132").Value = 2835.25
$ws.Range("J132").Value = 8888.857
$ws.Range("K132").Value = 8505.75
$ws.Range("L132").Value = 26666.571
$ws.Range("M132").Value = -5975.75
$ws.Range("N132").Value = -31726.571

$ws.Range("H136").Value = 13164247
$ws.Range("I136").Value = 15156668
$ws.Range("K136").Value = 45470004
$ws.Range("M136").Value = -45467454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3965.8462
$ws.Range("I3").Value = 999.3333
$ws.Range("J3").Value = 6508.5713
$ws.Range("K3").Value = 999.3333
$ws.Range("L3").Value = 6508.5713
$ws.Range("M3").Value = -885.3333
$ws.Range("N3").Value = -6736.5713

$ws.Range("H11").Value = 377.5
$ws.Range("I11").Value = 2.5
$ws.Range("J11").Value = 752.5
$ws.Range("K11").Value = 2.5
$ws.Range("L11").Value = 752.5
$ws.Range("M11").Value = 137.5
$ws.Range("N11").Value = -1032.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3632.9546
$ws.Range("J16").Value = 5161.6665
$ws.Range("L16").Value = 5161.6665
$ws.Range("N16").Value = -5735.6665

$ws.Range("H31").Value = 3910.6487
$ws.Range("I31").Value = 3320.8
$ws.Range("K31").Value = 3320.8
$ws.Range("M31").Value = -3025.8

$ws.Range("H34").Value = 3910.6487
$ws.Range("I34").Value = 3320.8
$ws.Range("K34").Value = 3320.8
$ws.Range("M34").Value = -3118.8

$ws.Range("H86").Value = 11617.571
$ws.Range("J86").Value = 10999
$ws.Range("L86").Value = 10999
$ws.Range("N86").Value = -13245

$ws.Range("H89").Value = 11617.571
$ws.Range("J89").Value = 10999
$ws.Range("L89").Value = 54995
$ws.Range("N89").Value = -66227

$ws.Range("H113").Value = 3632.9546
$ws.Range("J113").Value = 5161.6665
$ws.Range("L113").Value = 5161.6665
$ws.Range("N113").Value = -9501.666499999999

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H132").Value = 2778.1333
$ws.Range("I132").Value = 2833.7144
$ws.Range("K132").Value = 8501.143199999999
$ws.Range("M132").Value = -5971.143199999999

$ws.Range("H134").Value = 5727
$ws.Range("I134").Value = 5083.615
$ws.Range("J134").Value = 6563.4
$ws.Range("K134").Value = 15250.845
$ws.Range("L134").Value = 19690.2
$ws.Range("M134").Value = -12715.845
$ws.Range("N134").Value = -24760.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2587.25
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2587.25
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 7761.75
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -11601.75

$ws.Range("H128").Value = 509720
$ws.Range("I128").Value = 509720
$ws.Range("K128").Value = 1529160
$ws.Range("M128").Value = -1524180

$ws.Range("H131").Value = 23812752
$ws.Range("J131").Value = 4690
$ws.Range("L131").Value = 14070
$ws.Range("N131").Value = -24150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2475.889
$ws.Range("I102").Value = 1879.125
$ws.Range("K102").Value = 1879.125
$ws.Range("M102").Value = -257.125

$ws.Range("H126").Value = 4508.7
$ws.Range("I126").Value = 3514.5
$ws.Range("K126").Value = 10543.5
$ws.Range("M126").Value = -8073.5

$ws.Range("H132").Value = 3301.6785
$ws.Range("I132").Value = 2844.2273
$ws.Range("K132").Value = 8532.6819
$ws.Range("M132").Value = -6002.6819

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 55000
$ws.Range("J98").Value = 55000
$ws.Range("L98").Value = 55000
$ws.Range("N98").Value = -60990

$ws.Range("H101").Value = 61979.5
$ws.Range("J101").Value = 61979.5
$ws.Range("L101").Value = 61979.5
$ws.Range("N101").Value = -68469.5

$ws.Range("H128").Value = 99279.86
$ws.Range("J128").Value = 99279.86
$ws.Range("L128").Value = 99279.86
$ws.Range("N128").Value = -109239.86

$ws.Range("H132").Value = 4432.96
$ws.Range("I132").Value = 4335.452
$ws.Range("K132").Value = 13006.356
$ws.Range("M132").Value = -10476.356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 59958
$ws.Range("J63").Value = 59958
$ws.Range("L63").Value = 59958
$ws.Range("N63").Value = -61206

$ws.Range("H66").Value = 59958
$ws.Range("J66").Value = 59958
$ws.Range("L66").Value = 179874
$ws.Range("N66").Value = -186114

$ws.Range("H95").Value = 44999.5
$ws.Range("J95").Value = 44999.5
$ws.Range("L95").Value = 44999.5
$ws.Range("N95").Value = -50491.5

$ws.Range("H132").Value = 3278.8276
$ws.Range("I132").Value = 3403.625
$ws.Range("K132").Value = 10210.875
$ws.Range("M132").Value = -7680.875
